$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: insert "result_master_exclusion" before the existing "result_range_a" row (row 47) ---
$ws.Rows.Item(47).Insert()
$ws.Range("A47").Value = "result_master_exclusion"

# --- Block 2: insert 13 new "data_processed" codebook variable rows before "gender" (shifted to row 63) ---
$block2 = @(
    "bfi_a_mean",
    "bfi_c_mean",
    "bfi_e_mean",
    "bfi_n_mean",
    "bfi_o_mean",
    "exclude_bfi_completeness",
    "exclude_bfi_logical_scores",
    "exclude_iat_completeness",
    "exclude_iat_performance",
    "exclude_participant",
    "exclude_participant_without_bfi",
    "exclude_participant_without_iat",
    "exclude_unique_id_is.na"
)
$startRow = 63
for ($i = 0; $i -lt $block2.Length; $i++) {
    $r = $startRow + $i
    $ws.Rows.Item($r).Insert()
    $ws.Range("A$r").Value = $block2[$i]
}

# --- Block 3: insert 4 new IAT-summary codebook variable rows before "unique_id" (shifted to row 77) ---
$block3 = @(
    "iat_D",
    "iat_mean1",
    "iat_mean2",
    "iat_SD"
)
$startRow = 77
for ($i = 0; $i -lt $block3.Length; $i++) {
    $r = $startRow + $i
    $ws.Rows.Item($r).Insert()
    $ws.Range("A$r").Value = $block3[$i]
}
